$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = 5.997700000000001
    14  = 5.903000000000002
    16  = 5.334999999999999
    21  = 9.033799999999999
    23  = 8.903099999999997
    25  = 5.491000000000001
    26  = 5.038000000000006
    29  = 5.3006
    40  = 8.690800000000001
    53  = 5.704500000000001
    57  = 4.967699999999994
    59  = 5.229000000000001
    65  = 6.002499999999999
    69  = 5.438399999999994
    79  = 9.475700000000007
    83  = 5.387599999999995
    91  = 6.003400000000002
    93  = 5.658299999999997
    100 = 5.234400000000005
    103 = 5.858799999999998
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
